$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 582, shifting existing rows 582:609 down to 583:610
$ws.Rows.Item(582).Insert()

# Populate the newly inserted row 582 with the new weekly price record
$ws.Cells.Item(582, 1).Value = 5
$ws.Cells.Item(582, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(582, 3).Value = "Maule"
$ws.Cells.Item(582, 4).Value = 45147
$ws.Cells.Item(582, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(582, 5).Value = 7
$ws.Cells.Item(582, 6).Value = 100112032
$ws.Cells.Item(582, 7).Value = "Zapallo italiano"
$ws.Cells.Item(582, 8).Value = "Sin especificar"
$ws.Cells.Item(582, 9).Value = "Primera"
$ws.Cells.Item(582, 10).Value = 300
$ws.Cells.Item(582, 11).Value = 15000
$ws.Cells.Item(582, 12).Value = 15000
$ws.Cells.Item(582, 13).Value = 15000
$ws.Cells.Item(582, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(582, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(582, 16).Value = 300
$ws.Cells.Item(582, 17).Value = 50
$ws.Cells.Item(582, 18).Value = "Hortaliza"
